{"js": "// The edit removes the standalone \"\u0399.\u039a.\u03a5.\" bullet paragraph from the\n// \"\u039a\u039f\u0399\u039d\u039f\u03a0\u039f\u0399\u0397\u03a3\u0397\" (distribution) list, merging it away so the following\n// \"${local_directorate}\" bullet item takes its place in the list.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find((p) => p.text.trim() === \"\u0399.\u039a.\u03a5.\");\n\nif (target) {\n  // Paragraph.delete() removes the paragraph together with its paragraph\n  // mark, which merges it with the following paragraph \u2014 exactly what the\n  // diff shows (the \"\u0399.\u039a.\u03a5.\" paragraph and the following paragraph's own\n  // paragraph mark/properties both disappear, leaving a single paragraph\n  // that keeps this one's list formatting but the next one's run content).\n  target.delete();\n  await context.sync();\n}\n", "ps1": "# The edit removes the standalone \"\u0399.\u039a.\u03a5.\" bullet paragraph from the\n# \"\u039a\u039f\u0399\u039d\u039f\u03a0\u039f\u0399\u0397\u03a3\u0397\" (distribution) list. Deleting the paragraph's Range (text +\n# paragraph mark) merges it with the following paragraph, so the\n# \"${local_directorate}\" bullet item that followed it now takes its place\n# in the list - matching the target diff exactly.\n$d = $word.ActiveDocument\n\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute(\"\u0399.\u039a.\u03a5.\")\n\nif ($found) {\n    $paraRange = $d.Range($searchRange.Start, $searchRange.End)\n    $paraRange.Expand(4) | Out-Null   # wdParagraph\n    $paraRange.Delete()\n}\n"}
